$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format so values like "1.011" or "104.00"
# are stored as literal text (matching the source data) instead of being
# auto-converted to numbers by Excel's smart input parsing.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.210.89"
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("D3").Value = "1.835.82"
$ws.Range("E3").Value = "  +1.18%  "
$ws.Range("E4").Value = "  +1.12%  "
$ws.Range("D5").Value = "313.87"
$ws.Range("E5").Value = "  +1.37%  "
$ws.Range("D6").Value = "1.011"
$ws.Range("E6").Value = "  +1.05%  "
$ws.Range("D7").Value = "0.4713"
$ws.Range("E7").Value = "  +0.71%  "
$ws.Range("D8").Value = "0.3684"
$ws.Range("E8").Value = "  -0.28%  "
$ws.Range("D9").Value = "0.07428"
$ws.Range("E9").Value = "  +0.76%  "
$ws.Range("D10").Value = "0.8831"
$ws.Range("E10").Value = "  +1.49%  "
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("D12").Value = "1.814.11"
$ws.Range("E12").Value = "  +0.80%  "
$ws.Range("D13").Value = "0.07339"
$ws.Range("D14").Value = "5.481"
$ws.Range("E14").Value = "  +2.20%  "
$ws.Range("D15").Value = "92.96"
$ws.Range("E15").Value = "  +0.72%  "
$ws.Range("D16").Value = "6.574"
$ws.Range("E16").Value = "  +1.14%  "
$ws.Range("D17").Value = "1.013"
$ws.Range("E17").Value = "  +1.08%  "
$ws.Range("D18").Value = "0.000008797"
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("D19").Value = "1.011"
$ws.Range("E19").Value = "  +0.98%  "
$ws.Range("D20").Value = "14.80"
$ws.Range("E20").Value = "  +0.45%  "
$ws.Range("D21").Value = "27.229.38"
$ws.Range("E21").Value = "  +1.25%  "
$ws.Range("D22").Value = "5.312"
$ws.Range("E22").Value = "  -0.46%  "
$ws.Range("D23").Value = "10.70"
$ws.Range("E23").Value = "  +1.56%  "
$ws.Range("D24").Value = "2.057.45"
$ws.Range("E24").Value = "  +0.86%  "
$ws.Range("E25").Value = "  +0.61%  "
$ws.Range("D26").Value = "152.76"
$ws.Range("E26").Value = "  +0.64%  "
$ws.Range("D27").Value = "18.63"
$ws.Range("E27").Value = "  +1.56%  "
$ws.Range("D28").Value = "2.174"
$ws.Range("E28").Value = "  -0.77%  "
$ws.Range("D29").Value = "5.283"
$ws.Range("E29").Value = "  -0.23%  "
$ws.Range("D30").Value = "117.81"
$ws.Range("E30").Value = "  +2.02%  "
$ws.Range("D31").Value = "0.08929"
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("D32").Value = "0.7608"
$ws.Range("E32").Value = "  -0.53%  "
$ws.Range("D33").Value = "1.173"
$ws.Range("E33").Value = "  +1.11%  "
$ws.Range("D34").Value = "4.546"
$ws.Range("E34").Value = "  +1.50%  "
$ws.Range("E35").Value = "  +0.85%  "
$ws.Range("E36").Value = "  +1.08%  "
$ws.Range("D37").Value = "1.106"
$ws.Range("E37").Value = "  +1.21%  "
$ws.Range("D38").Value = "0.05342"
$ws.Range("E38").Value = "  +1.42%  "
$ws.Range("D39").Value = "0.01961"
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("D40").Value = "3.004"
$ws.Range("E40").Value = "  +2.03%  "
$ws.Range("D41").Value = "7.335"
$ws.Range("E41").Value = "  +1.35%  "
$ws.Range("D42").Value = "2.405"
$ws.Range("E42").Value = "  +1.94%  "
$ws.Range("D43").Value = "0.5355"
$ws.Range("E43").Value = "  +0.40%  "
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("D45").Value = "8.553"
$ws.Range("E45").Value = "  +1.27%  "
$ws.Range("D46").Value = "0.4959"
$ws.Range("E46").Value = "  +0.45%  "
$ws.Range("D47").Value = "10.53"
$ws.Range("E47").Value = "  +0.83%  "
$ws.Range("D48").Value = "1.012"
$ws.Range("E48").Value = "  +1.15%  "
$ws.Range("D49").Value = "1.674"
$ws.Range("E49").Value = "  +0.31%  "
$ws.Range("D50").Value = "104.04"
$ws.Range("E50").Value = "  +1.28%  "
$ws.Range("D51").Value = "0.06326"
$ws.Range("E51").Value = "  +0.71%  "

# Restore the default (unstyled) cell format on column D now that the
# text values are safely stored, so the style table matches the original.
$ws.Range("D2:D51").Style = "Normal"

